$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link text updates (columns B and C) ---
# Plain text values; Excel never mistakes these for numbers, so a simple
# Range.Value assignment is all that is needed.
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

# --- Volume(1h) percentage updates (column E) ---
# Values are padded with spaces and a trailing "%", so Excel keeps them as
# text automatically.
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").Value = "  +6.42%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("E38").Value = "  +4.31%  "
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("E42").Value = "  +7.40%  "
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("E51").Value = "  -0.24%  "

# --- Price updates (column D) ---
# Several of these look like plain numbers (e.g. "10.20", "6.50"), and a bare
# Range.Value assignment would let Excel "helpfully" reinterpret them as
# numeric values, silently dropping significant trailing zeros (10.20 -> 10.2)
# and turning the cell from Text into a Number type. The source workbook
# stores every Price cell as literal text, so for each touched Price cell we
# force Text format first, assign the literal string, then clear the
# temporary formatting again so the cell ends up with no explicit style --
# exactly like every other (untouched) Price cell in the sheet.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "59.155.76"
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.593.82"
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "521.82"
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "140.30"
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.612.27"
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "6.51"
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "3.058.41"
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "59.166.76"
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "20.40"
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.595.92"
$cell.ClearFormats()
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "338.68"
$cell.ClearFormats()
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "4.33"
$cell.ClearFormats()
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.20"
$cell.ClearFormats()
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.50"
$cell.ClearFormats()
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "66.40"
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.ClearFormats()
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.05"
$cell.ClearFormats()
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.ClearFormats()
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0727"
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "5.95"
$cell.ClearFormats()
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "18.84"
$cell.ClearFormats()
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "149.05"
$cell.ClearFormats()
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.00"
$cell.ClearFormats()
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "36.33"
$cell.ClearFormats()
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.46"
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.831"
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.825"
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "276.44"
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "10.73"
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0952"
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.589"
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "18.60"
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.982.83"
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "4.63"
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0220"
$cell.ClearFormats()
